$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2; existing rows 2-9 (transactions 1-8)
# shift down to rows 3-10.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the header's formatting. Re-apply the
# plain bordered data-row style (used by every other data row) by copying
# the format from the row now directly beneath it (the old row 2).
$ws.Range("A3:J3").Copy()
$ws.Range("A2:J2").PasteSpecial(-4122)

# Populate the new "Opening Balance" row.
$ws.Range("D2").Value = "Opening Balance"

# Leading apostrophe forces the numeric-looking balance to stay literal text
# (matching the "20000.00" text already used throughout this sheet).
$ws.Range("G2").Value = "'20000.00"

# Writing the quoted value nudges G2's style (quote-prefix flag); restore
# the plain data-row style by re-pasting formats from its neighbor.
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)
